$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 1 de Junio de 2020 a las 16:05"

# --- Update Estados Unidos (row 4) ---
$ws.Cells.Item(4, 2).Value = 1839126
$ws.Cells.Item(4, 3).Value = 1956
$ws.Cells.Item(4, 5).Value = 1133010
$ws.Cells.Item(4, 7).Value = 46
$ws.Cells.Item(4, 8).Value = 106241

# --- Update India (row 10) ---
$ws.Cells.Item(10, 2).Value = 192782
$ws.Cells.Item(10, 3).Value = 2173
$ws.Cells.Item(10, 4).Value = 92474
$ws.Cells.Item(10, 5).Value = 94893

# --- Update Alemania (row 12) ---
$ws.Cells.Item(12, 2).Value = 183564
$ws.Cells.Item(12, 3).Value = 70
$ws.Cells.Item(12, 5).Value = 9059

# --- Sudafrica overtakes Portugal: swap rows 31/32 ---
# Row 31 becomes Sudafrica with brand-new figures
$ws.Cells.Item(31, 1).Value = "Sudafrica"
$ws.Cells.Item(31, 2).Value = 34357
$ws.Cells.Item(31, 3).Value = 1674
$ws.Cells.Item(31, 4).Value = 16809
$ws.Cells.Item(31, 5).Value = 16843
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(31, 7).Value = 22
$ws.Cells.Item(31, 8).Value = 705

# Row 32 becomes Portugal, keeping its previous figures
$ws.Cells.Item(32, 1).Value = "Portugal"
$ws.Cells.Item(32, 2).Value = 32700
$ws.Cells.Item(32, 3).Value = 200
$ws.Cells.Item(32, 4).Value = 19552
$ws.Cells.Item(32, 5).Value = 11724
$ws.Cells.Item(32, 6).Value = 0
$ws.Cells.Item(32, 7).Value = 14
$ws.Cells.Item(32, 8).Value = 1424

# --- Rumania (row 41) ---
$ws.Cells.Item(41, 5).Value = 4696
$ws.Cells.Item(41, 7).Value = 10
$ws.Cells.Item(41, 8).Value = 1276

# --- Serbia (row 54) ---
$ws.Cells.Item(54, 2).Value = 11430
$ws.Cells.Item(54, 3).Value = 18
$ws.Cells.Item(54, 4).Value = 6726
$ws.Cells.Item(54, 5).Value = 4460
$ws.Cells.Item(54, 7).Value = 1
$ws.Cells.Item(54, 8).Value = 244

# --- Row 61 ---
$ws.Cells.Item(61, 2).Value = 8442
$ws.Cells.Item(61, 3).Value = 2
$ws.Cells.Item(61, 5).Value = 479

# --- Row 62 (Moldavia) ---
$ws.Cells.Item(62, 4).Value = 4622
$ws.Cells.Item(62, 5).Value = 3331
$ws.Cells.Item(62, 7).Value = 3
$ws.Cells.Item(62, 8).Value = 298

# --- Row 67 (Finlandia) ---
$ws.Cells.Item(67, 5).Value = 1067
$ws.Cells.Item(67, 8).Value = 318

# --- Row 69 (Camerun) ---
$ws.Cells.Item(69, 2).Value = 6143
$ws.Cells.Item(69, 3).Value = 239
$ws.Cells.Item(69, 4).Value = 3578
$ws.Cells.Item(69, 5).Value = 2368
$ws.Cells.Item(69, 7).Value = 6
$ws.Cells.Item(69, 8).Value = 197

# --- Row 75 (Tayikistan) ---
$ws.Cells.Item(75, 2).Value = 4013
$ws.Cells.Item(75, 3).Value = 83
$ws.Cells.Item(75, 4).Value = 2089
$ws.Cells.Item(75, 5).Value = 1877

# --- Row 80 (Republica de Yibuti) ---
$ws.Cells.Item(80, 2).Value = 3569
$ws.Cells.Item(80, 3).Value = 215
$ws.Cells.Item(80, 4).Value = 1521
$ws.Cells.Item(80, 5).Value = 2024

# --- Nepal overtakes Islandia: swap rows 97/98 ---
# Row 97 becomes Nepal with brand-new figures
$ws.Cells.Item(97, 1).Value = "Nepal"
$ws.Cells.Item(97, 2).Value = 1811
$ws.Cells.Item(97, 3).Value = 239
$ws.Cells.Item(97, 4).Value = 221
$ws.Cells.Item(97, 5).Value = 1582
$ws.Cells.Item(97, 6).Value = 0
$ws.Cells.Item(97, 7).Value = 0
$ws.Cells.Item(97, 8).Value = 8

# Row 98 becomes Islandia, keeping its previous figures
$ws.Cells.Item(98, 1).Value = "Islandia"
$ws.Cells.Item(98, 2).Value = 1806
$ws.Cells.Item(98, 3).Value = 0
$ws.Cells.Item(98, 4).Value = 1794
$ws.Cells.Item(98, 5).Value = 2
$ws.Cells.Item(98, 6).Value = 0
$ws.Cells.Item(98, 7).Value = 0
$ws.Cells.Item(98, 8).Value = 10

# --- Row 102 (Sri Lanka) ---
$ws.Cells.Item(102, 2).Value = 1639
$ws.Cells.Item(102, 3).Value = 6
$ws.Cells.Item(102, 5).Value = 817

# --- Row 148 ---
$ws.Cells.Item(148, 4).Value = 311
$ws.Cells.Item(148, 5).Value = 1
